# Updated the logic_BOM to have the SSR Board

$wb = $excel.ActiveWorkbook
$wsLogic = $wb.Worksheets.Item(1)
$wsSSR = $wb.Worksheets.Item(2)

$wsLogic.Name = "Logic"
$wsSSR.Name = "SSR"

# ---- SSR sheet: header row ----
$wsSSR.Range("A1").Value = "Quantity"
$wsSSR.Range("B1").Value = "Value"
$wsSSR.Range("C1").Value = "Device"
$wsSSR.Range("D1").Value = "Package"
$wsSSR.Range("E1").Value = "DigiKey ID"

# ---- SSR sheet: Quantity column ----
$wsSSR.Range("A2").Value = 1
$wsSSR.Range("A3").Value = 1
$wsSSR.Range("A4").Value = 1
$wsSSR.Range("A5").Value = 1
$wsSSR.Range("A6").Value = 1
$wsSSR.Range("A7").Value = 1
$wsSSR.Range("A8").Value = 1
$wsSSR.Range("A9").Value = 1
$wsSSR.Range("A10").Value = 1
$wsSSR.Range("A11").Value = 1
$wsSSR.Range("A12").Value = 1

# ---- SSR sheet: parts list (entered in documentation order) ----
$wsSSR.Range("B2").Value = "acs770"
$wsSSR.Range("C3").Value = "High Power Terminal"
$wsSSR.Range("C2").Value = "current sense"
$wsSSR.Range("B4").Value = "D_STP5045S"
$wsSSR.Range("C4").Value = "Diode"
$wsSSR.Range("D3").Value = "Through hole"
$wsSSR.Range("B5").Value = "60k"
$wsSSR.Range("C5").Value = "Resistor"
$wsSSR.Range("B7").Value = "SPDT"
$wsSSR.Range("C7").Value = "Switch"
$wsSSR.Range("B8").Value = "PTSCMD"
$wsSSR.Range("C8").Value = "Poly-fuse"
$wsSSR.Range("B9").Value = "lt1910"
$wsSSR.Range("C9").Value = "Mosfet driver"
$wsSSR.Range("D9").Value = "SOIC-8"
$wsSSR.Range("B10").Value = "irls3034-7p"
$wsSSR.Range("C10").Value = "Mosfet"
$wsSSR.Range("D10").Value = "Surface mount"
$wsSSR.Range("B11").Value = "10uF"
$wsSSR.Range("C11").Value = "Capacitor"
$wsSSR.Range("E12").Value = "H3373-ND"
$wsSSR.Range("C12").Value = "Connector"
$wsSSR.Range("F12").Value = "FIND THE MALE CONNECTOR"

# ---- SSR sheet: remaining cells (duplicate text / numeric package sizes) ----
$wsSSR.Range("D2").Value = "Surface mount"
$wsSSR.Range("D4").Value = "Through hole"
$wsSSR.Range("D5").Value = 1206
$wsSSR.Range("B6").Value = "10k"
$wsSSR.Range("C6").Value = "Resistor"
$wsSSR.Range("D6").Value = 1206
$wsSSR.Range("D7").Value = "Through hole"
$wsSSR.Range("D8").Value = 1206
$wsSSR.Range("D11").Value = 1206
$wsSSR.Range("D12").Value = "Surface mount"

# ---- SSR sheet: column widths (best-fit) ----
$wsSSR.Columns.Item(3).ColumnWidth = 18.85
$wsSSR.Columns.Item(4).ColumnWidth = 13.1667
$wsSSR.Columns.Item(6).ColumnWidth = 25.6

# ---- SSR sheet: hyperlink for the DigiKey part ----
$digikeyUrl = "http://www.digikey.com/scripts/DkSearch/dksus.dll?Detail&itemSeq=165378352&uq=635580674242796387&CSRT=10282226202437667372"
$wsSSR.Hyperlinks.Add($wsSSR.Range("E12"), $digikeyUrl, [Type]::Missing, [Type]::Missing, $digikeyUrl)
$wsSSR.Range("E12").Value = "H3373-ND"

# ---- Logic sheet: selection moves to A1:E1, no longer the active tab ----
$wsLogic.Activate()
$wsLogic.Range("A1:E1").Select() | Out-Null

# ---- SSR sheet: becomes the active tab with I15 selected ----
$wsSSR.Activate()
$wsSSR.Range("I15").Select() | Out-Null
